# Scheduled runner update: refresh cached market-board price/profit figures
# across the Zodiark_Profits crafting-leve sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Writes the recomputed currentAveragePrice* / LevePrice* /
# LeveProfit* cached numbers back into each sheet's Table_<JOB> range.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1542.0164
$ws.Range("J17").Value = 1583.9828
$ws.Range("L17").Value = 4751.9484
$ws.Range("N17").Value = -5087.9484
$ws.Range("H28").Value = 2246.1667
$ws.Range("I28").Value = 2174
$ws.Range("J28").Value = 2498.75
$ws.Range("K28").Value = 2174
$ws.Range("L28").Value = 2498.75
$ws.Range("M28").Value = -1689
$ws.Range("N28").Value = -3468.75
$ws.Range("H32").Value = 4287.8945
$ws.Range("I32").Value = 3237.6667
$ws.Range("J32").Value = 4772.615
$ws.Range("K32").Value = 3237.6667
$ws.Range("L32").Value = 4772.615
$ws.Range("M32").Value = -2911.6667
$ws.Range("N32").Value = -5424.615
$ws.Range("H64").Value = 7717
$ws.Range("J64").Value = 8645.556
$ws.Range("L64").Value = 8645.556
$ws.Range("N64").Value = -9141.556
$ws.Range("H67").Value = 7717
$ws.Range("J67").Value = 8645.556
$ws.Range("L67").Value = 8645.556
$ws.Range("N67").Value = -10361.556
$ws.Range("H97").Value = 2802.6428
$ws.Range("J97").Value = 2956.7693
$ws.Range("L97").Value = 8870.3079
$ws.Range("N97").Value = -9862.3079
$ws.Range("H98").Value = 1405.4166
$ws.Range("I98").Value = 1421.55
$ws.Range("J98").Value = 1324.75
$ws.Range("K98").Value = 1421.55
$ws.Range("L98").Value = 1324.75
$ws.Range("M98").Value = 76.45000000000005
$ws.Range("N98").Value = -4320.75
$ws.Range("H112").Value = 1285.38
$ws.Range("J112").Value = 1285.38
$ws.Range("L112").Value = 3856.14
$ws.Range("N112").Value = -6072.14
$ws.Range("H122").Value = 1405.4166
$ws.Range("I122").Value = 1421.55
$ws.Range("J122").Value = 1324.75
$ws.Range("K122").Value = 4264.65
$ws.Range("L122").Value = 3974.25
$ws.Range("M122").Value = -1814.65
$ws.Range("N122").Value = -8874.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 9999.5
$ws.Range("J92").Value = 9999.5
$ws.Range("L92").Value = 9999.5
$ws.Range("N92").Value = -14991.5
$ws.Range("H107").Value = 38614
$ws.Range("J107").Value = 38614
$ws.Range("L107").Value = 38614
$ws.Range("N107").Value = -46294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2395.3809
$ws.Range("I107").Value = 1657.4706
$ws.Range("J107").Value = 5531.5
$ws.Range("K107").Value = 1657.4706
$ws.Range("L107").Value = 5531.5
$ws.Range("M107").Value = 262.5293999999999
$ws.Range("N107").Value = -9371.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3649.0417
$ws.Range("I16").Value = 4232.1
$ws.Range("J16").Value = 733.75
$ws.Range("K16").Value = 4232.1
$ws.Range("L16").Value = 733.75
$ws.Range("M16").Value = -3945.1
$ws.Range("N16").Value = -1307.75
$ws.Range("H99").Value = 5162.9
$ws.Range("I99").Value = 1931.5
$ws.Range("K99").Value = 1931.5
$ws.Range("M99").Value = -433.5
$ws.Range("H113").Value = 3649.0417
$ws.Range("I113").Value = 4232.1
$ws.Range("J113").Value = 733.75
$ws.Range("K113").Value = 4232.1
$ws.Range("L113").Value = 733.75
$ws.Range("M113").Value = -2062.1
$ws.Range("N113").Value = -5073.75
$ws.Range("H122").Value = 2356.5715
$ws.Range("I122").Value = 1923.2307
$ws.Range("K122").Value = 5769.6921
$ws.Range("M122").Value = -3319.6921
$ws.Range("H126").Value = 5162.9
$ws.Range("I126").Value = 1931.5
$ws.Range("K126").Value = 5794.5
$ws.Range("M126").Value = -3324.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 472.25
$ws.Range("I33").Value = 340.57144
$ws.Range("J33").Value = 574.6667
$ws.Range("K33").Value = 2043.42864
$ws.Range("L33").Value = 3448.0002
$ws.Range("M33").Value = -1760.42864
$ws.Range("N33").Value = -4014.0002
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H138").Value = 7104
$ws.Range("I138").Value = 7104
$ws.Range("K138").Value = 21312
$ws.Range("M138").Value = -16172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H122").Value = 2994.647
$ws.Range("I122").Value = 1903.5834
$ws.Range("J122").Value = 5613.2
$ws.Range("K122").Value = 5710.7502
$ws.Range("L122").Value = 16839.6
$ws.Range("M122").Value = -3260.7502
$ws.Range("N122").Value = -21739.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 52504.5
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H61").Value = 1728.3158
$ws.Range("I61").Value = 1797.8
$ws.Range("J61").Value = 1467.75
$ws.Range("K61").Value = 1797.8
$ws.Range("L61").Value = 1467.75
$ws.Range("M61").Value = -1595.8
$ws.Range("N61").Value = -1871.75
$ws.Range("H106").Value = 15898.5
$ws.Range("J106").Value = 15898.5
$ws.Range("L106").Value = 15898.5
$ws.Range("N106").Value = -18422.5
$ws.Range("H113").Value = 1728.3158
$ws.Range("I113").Value = 1797.8
$ws.Range("J113").Value = 1467.75
$ws.Range("K113").Value = 1797.8
$ws.Range("L113").Value = 1467.75
$ws.Range("M113").Value = 372.2
$ws.Range("N113").Value = -5807.75
$ws.Range("H122").Value = 5495.875
$ws.Range("I122").Value = 4773.077
$ws.Range("J122").Value = 6350.091
$ws.Range("K122").Value = 14319.231
$ws.Range("L122").Value = 19050.273
$ws.Range("M122").Value = -11869.231
$ws.Range("N122").Value = -23950.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4087.2
$ws.Range("I81").Value = 1904.7
$ws.Range("J81").Value = 8452.200000000001
$ws.Range("K81").Value = 3809.4
$ws.Range("L81").Value = 16904.4
$ws.Range("M81").Value = -2748.4
$ws.Range("N81").Value = -19026.4
$ws.Range("H84").Value = 4087.2
$ws.Range("I84").Value = 1904.7
$ws.Range("J84").Value = 8452.200000000001
$ws.Range("K84").Value = 19047
$ws.Range("L84").Value = 84522
$ws.Range("M84").Value = -13743
$ws.Range("N84").Value = -95130
$ws.Range("H126").Value = 22223922
$ws.Range("I126").Value = 23811202
$ws.Range("J126").Value = 1998
$ws.Range("K126").Value = 71433606
$ws.Range("L126").Value = 5994
$ws.Range("M126").Value = -71431136
$ws.Range("N126").Value = -10934
